$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata": update Date value (B8) ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2025-07-21T11:52:46+00:00"

# --- Sheet "Include #0": update System URI value (B4) ---
$wsInc0 = $wb.Worksheets.Item("Include #0")
$wsInc0.Range("B4").Value = "https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R97-DroitExerciceCompl"

# --- Sheet "Include #1": update System URI value (B4) ---
$wsInc1 = $wb.Worksheets.Item("Include #1")
$wsInc1.Range("B4").Value = "https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R01-EnsembleSavoirFaire-CISIS"
